$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "-"

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "-"

# Row 4
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"

# Row 6
$ws.Range("C6").Value = "-"
$ws.Range("E6").Value = "-"

# Row 7
$ws.Range("C7").Value = "-"
$ws.Range("E7").Value = "-"

# Row 11
$ws.Range("E11").Value = "[Aline S. M.-T. M. Metalicos-1A, -, -, -]"

# Row 12
$ws.Range("B12").Value = "Anselmo-Gestao Integrada"
$ws.Range("E12").Value = "[Aline S. M.-T. M. Metalicos-1A, Andre B.-Comandos Eletricos-1A, Andre B.-Comandos Eletricos-1A, Andre B.-Comandos Eletricos-1A]"

# Row 14
$ws.Range("B14").Value = "Anselmo-Gestao Integrada"
$ws.Range("E14").Value = "[Ismail-Metrologia 1-1A, Ismail-Metrologia 1-1A, Ismail-Metrologia 1-1A, Andre B.-Comandos Eletricos-1A]"

# Row 15
$ws.Range("E15").Value = "[Elcio D.-Desenho tecnico mecanico-1A, Ismail-Metrologia 1-1A, Elcio D.-Desenho tecnico mecanico-1A, Elcio D.-Desenho tecnico mecanico-1A]"

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "-"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = "-"

# Row 20
$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "-"
